$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '43.359.49'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -1.71%  '
# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.250.41'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -0.61%  '
# Row 4
$ws.Range("E4").Value = '  +0.19%  '
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '230.74'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.21%  '
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.639'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +1.12%  '
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '64.23'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +1.02%  '
# Row 8
$ws.Range("E8").Value = '  +0.07%  '
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.439'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -1.64%  '
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0953'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -8.14%  '
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '56.78'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -0.15%  '
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '26.59'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +1.05%  '
# Row 13
$ws.Range("E13").Value = '  -1.56%  '
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.584.54'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -0.59%  '
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '14.91'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -5.14%  '
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.02'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -0.61%  '
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.821'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -2.05%  '
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.246.55'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -1.12%  '
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '43.289.25'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -1.43%  '
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0964'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -4.79%  '
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '73.06'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -0.75%  '
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.04'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +0.15%  '
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '246.58'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -3.20%  '
# Row 24
$ws.Range("E24").Value = '  -0.02%  '
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.79'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +13.36%  '
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.41'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -0.74%  '
# Row 27
$ws.Range("E27").Value = '  -1.58%  '
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.70'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -4.34%  '
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '173.76'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +1.09%  '
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '21.56'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +3.46%  '
# Row 31
$ws.Range("B31").Value = 'ImmutableX'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.43'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +2.73%  '
# Row 32
$ws.Range("B32").Value = 'Kaspa'
$ws.Range("C32").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.130'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -5.51%  '
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.124'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +0.70%  '
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.92'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +3.08%  '
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0677'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -0.59%  '
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.92'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +1.13%  '
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.65'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -5.17%  '
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.36'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -5.29%  '
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.26'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -2.49%  '
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0249'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -2.75%  '
# Row 41
$ws.Range("E41").Value = '  -0.23%  '
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.76'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +5.17%  '
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.49'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +2.58%  '
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '17.24'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -0.88%  '
# Row 45
$ws.Range("E45").Value = '  +1.86%  '
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '96.57'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -1.14%  '
# Row 47
$ws.Range("E47").Value = '  -0.82%  '
# Row 48
$ws.Range("E48").Value = '  -3.02%  '
# Row 49
$ws.Range("E49").Value = '  -2.45%  '
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.426.81'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -1.42%  '
# Row 51
$ws.Range("B51").Value = 'HuobiToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.74'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +0.20%  '
